$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.908.59"
$ws.Cells.Item(2, 5).Value = "  +0.34%  "
$ws.Cells.Item(3, 4).Value = "2.440.13"
$ws.Cells.Item(3, 5).Value = "  -0.76%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Value = "'560.46"
$ws.Cells.Item(5, 5).Value = "  +0.24%  "
$ws.Cells.Item(6, 4).Value = "'162.14"
$ws.Cells.Item(6, 5).Value = "  +0.31%  "
$ws.Cells.Item(7, 5).Value = "  -0.03%  "
$ws.Cells.Item(8, 5).Value = "  +1.98%  "
$ws.Cells.Item(9, 5).Value = "  +11.34%  "
$ws.Cells.Item(10, 5).Value = "  -1.61%  "
$ws.Cells.Item(11, 5).Value = "  -0.02%  "
$ws.Cells.Item(12, 5).Value = "  -5.07%  "
$ws.Cells.Item(13, 5).Value = "  +4.85%  "
$ws.Cells.Item(14, 4).Value = "68.789.26"
$ws.Cells.Item(14, 5).Value = "  +0.28%  "
$ws.Cells.Item(15, 4).Value = "2.876.70"
$ws.Cells.Item(15, 5).Value = "  -1.00%  "
$ws.Cells.Item(16, 4).Value = "'23.18"
$ws.Cells.Item(16, 5).Value = "  -1.37%  "
$ws.Cells.Item(17, 4).Value = "2.427.15"
$ws.Cells.Item(17, 5).Value = "  -0.62%  "
$ws.Cells.Item(18, 4).Value = "'10.50"
$ws.Cells.Item(18, 5).Value = "  -1.26%  "
$ws.Cells.Item(19, 4).Value = "'338.64"
$ws.Cells.Item(19, 5).Value = "  +1.36%  "
$ws.Cells.Item(20, 4).Value = "'6.94"
$ws.Cells.Item(20, 5).Value = "  +0.61%  "
$ws.Cells.Item(21, 4).Value = "'3.83"
$ws.Cells.Item(21, 5).Value = "  +1.40%  "
$ws.Cells.Item(22, 5).Value = "  +2.75%  "
$ws.Cells.Item(23, 5).Value = "  +0.01%  "
$ws.Cells.Item(24, 4).Value = "'67.04"
$ws.Cells.Item(24, 5).Value = "  +0.89%  "
$ws.Cells.Item(25, 5).Value = "  +1.90%  "
$ws.Cells.Item(26, 4).Value = "2.553.83"
$ws.Cells.Item(26, 5).Value = "  -1.30%  "
$ws.Cells.Item(27, 4).Value = "'1.00"
$ws.Cells.Item(27, 5).Value = "  -0.05%  "
$ws.Cells.Item(28, 4).Value = "'8.22"
$ws.Cells.Item(28, 5).Value = "  +0.94%  "
$ws.Cells.Item(29, 4).Value = "0.0₃0817"
$ws.Cells.Item(29, 5).Value = "  +0.59%  "
$ws.Cells.Item(30, 4).Value = "'7.13"
$ws.Cells.Item(30, 5).Value = "  -0.41%  "
$ws.Cells.Item(31, 4).Value = "'1.00"
$ws.Cells.Item(31, 5).Value = "  +0.04%  "
$ws.Cells.Item(32, 4).Value = "'428.28"
$ws.Cells.Item(32, 5).Value = "  -0.15%  "
$ws.Cells.Item(33, 5).Value = "  +1.97%  "
$ws.Cells.Item(34, 5).Value = "  +0.36%  "
$ws.Cells.Item(35, 4).Value = "'159.56"
$ws.Cells.Item(35, 5).Value = "  +0.38%  "
$ws.Cells.Item(36, 5).Value = "  +0.10%  "
$ws.Cells.Item(37, 5).Value = "  +0.01%  "
$ws.Cells.Item(38, 4).Value = "'17.98"
$ws.Cells.Item(38, 5).Value = "  +1.48%  "
$ws.Cells.Item(39, 5).Value = "  -2.84%  "
$ws.Cells.Item(40, 4).Value = "'0.297"
$ws.Cells.Item(40, 5).Value = "  -0.71%  "
$ws.Cells.Item(41, 5).Value = "  +3.71%  "
$ws.Cells.Item(42, 4).Value = "'4.34"
$ws.Cells.Item(42, 5).Value = "  -1.36%  "
$ws.Cells.Item(43, 5).Value = "  +0.97%  "
$ws.Cells.Item(44, 5).Value = "  -0.84%  "
$ws.Cells.Item(45, 4).Value = "'3.35"
$ws.Cells.Item(45, 5).Value = "  +0.20%  "
$ws.Cells.Item(46, 4).Value = "'130.56"
$ws.Cells.Item(46, 5).Value = "  +0.53%  "
$ws.Cells.Item(47, 4).Value = "'0.0717"
$ws.Cells.Item(47, 5).Value = "  +0.64%  "
$ws.Cells.Item(48, 5).Value = "  +0.05%  "
$ws.Cells.Item(49, 5).Value = "  -0.15%  "
$ws.Cells.Item(50, 4).Value = "'0.0923"
$ws.Cells.Item(50, 5).Value = "  +1.74%  "
$ws.Cells.Item(51, 5).Value = "  +1.15%  "
